$wb = $excel.ActiveWorkbook

# Work on the "Repayment schedule" worksheet
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column before column N (shifts old N/O/P -> O/P/Q)
$ws.Columns("N").Insert()
$ws.Columns("N").ColumnWidth = 10.1666666666667

# Activate the "Repayment schedule" sheet and select cell R6 to match the
# recorded selection/active tab state
$ws.Activate()
$ws.Range("R6").Select()
